$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TRUID values in the two "Project_NeonUser1" rows (A2, A3)
$ws.Range("I2").Value = "USER_NAME=Project_NeonUser1@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_FIRST_NAME=Project||USER_LAST_NAME=Neon1||USER_MIDDLE_NAME=TR||USER_TRUID=944694d8-f72a-4be4-906c-1e53d3232098"
$ws.Range("I3").Value = "USER_NAME=Project_NeonUser1@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_FIRST_NAME=Project||USER_LAST_NAME=Neon1||USER_MIDDLE_NAME=TR||USER_TRUID=444694d8-f72a-4be4-906c-1e53d3232098"

# Swap the JDRUser5 / JDRUser6 rows and update TRUID values
$ws.Range("I4").Value = "USER_NAME=Neon_JDRUser6@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_TRUID=555694d8-f72a-4be4-906c-1e53d3235466"
$ws.Range("I5").Value = "USER_NAME=Neon_JDRUser5@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_FIRST_NAME=JDR||USER_LAST_NAME=E5||USER_MIDDLE_NAME=REDDY5||USER_TRUID=555694d8-f72a-4be4-906c-1e53d3235466"

# Update the selected cell in the sheet view
$ws.Range("G3").Select()
